# Add the "CarDetails" worksheet (upcoming Tata car details) after the
# existing "PopularCarModels" sheet, matching the BikeDetails/PopularCarModels
# header styling (white text on a green fill).

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CarDetails"

# ---- header row -----------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Car Name"
$ws.Cells.Item(1, 2).Value = "Car Price"
$ws.Cells.Item(1, 3).Value = "Launch Date"

$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Color = 16777215
$headerRange.Interior.Color = 32768

# ---- data rows --------------------------------------------------------
# Cells that look like "Mon YYYY" (no day-of-month) get auto-parsed as
# dates by the host when assigned through .Value; prefix them with an
# apostrophe to force literal text, matching the source workbook where
# every cell in this sheet is a shared string.
$rows = @(
    @("Tata Altroz Racer",     "Rs. 10.00 Lakh", "20 Mar 2024"),
    @("Tata Curvv EV",         "Rs. 20.00 Lakh", "'Jul 2024"),
    @("Tata Curvv",            "Rs. 10.50 Lakh", "'Aug 2024"),
    @("Tata Avinya",           "Rs. 30.00 Lakh", "'Jan 2025"),
    @("Tata Harrier EV",       "Rs. 30.00 Lakh", "'Apr 2025"),
    @("Tata Punch 2025",       "Rs. 6.00 Lakh",  "'Jun 2025"),
    @("Tata Sierra",           "Rs. 25.00 Lakh", "'Dec 2025"),
    @("Tata Kite 5",           "Rs. 4.50 Lakh",  "Unrevealed"),
    @("Tata Atmos",            "Rs. 12.00 Lakh", "Unrevealed"),
    @("Tata H7X",              "Rs. 15.00 Lakh", "Unrevealed"),
    @("Tata Altroz EV",        "Rs. 14.00 Lakh", "Unrevealed"),
    @("Tata Hexa",             "Rs. 14.00 Lakh", "Unrevealed"),
    @("Tata EVision Electric", "Rs. 25.00 Lakh", "Unrevealed")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

# Strip the quote-prefix marker left behind by the apostrophe workaround
# above so the data rows end up on the plain default style, same as the
# source file.
$ws.Range("A2:C14").ClearFormats()

# ---- column widths ----------------------------------------------------
# Best-fit widths (in "characters") matching the source file's stored
# column widths of 17.0390625 / 12.37109375 / 11.140625 (256ths of a char).
$ws.Columns.Item(1).ColumnWidth = 16.09
$ws.Columns.Item(2).ColumnWidth = 11.42
$ws.Columns.Item(3).ColumnWidth = 10.25
